# Atualizando o arquivo XLSX
# Refresh today's fixtures/odds: the match that used to sit in row 2
# (Honduras Liga Nacional) is replaced by the Atlas x Monterrey friendly
# that previously occupied row 3, and row 3 now holds a new fixture
# (Alebrijes de Oaxaca x Chapulineros de Oaxaca) with its own odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range('A2').Value = 'Friendly Matches'
$ws.Range('C2').Value = '21:30:00'
$ws.Range('D2').Value = 'Atlas'
$ws.Range('E2').Value = 'Monterrey'
$ws.Range('F2').Value = 4.3
$ws.Range('G2').Value = 5.4
$ws.Range('H2').Value = 1.83
$ws.Range('I2').Value = 2
$ws.Range('J2').Value = 3.55
$ws.Range('K2').Value = 4.6
$ws.Range('L2').Value = 1.35
$ws.Range('M2').Value = 1.08
$ws.Range('N2').Value = 3.7
$ws.Range('O2').Value = 1.27
$ws.Range('P2').Value = 1.83
$ws.Range('Q2').Value = 2.12
$ws.Range('R2').Value = 1.32
$ws.Range('S2').Value = 4
$ws.Range('T2').Value = 1.77
$ws.Range('U2').Value = 1.75
$ws.Range('V2').Value = 2.02
$ws.Range('W2').Value = 1.23
$ws.Range('X2').Value = 12
$ws.Range('AB2').Value = 1000
$ws.Range('AC2').Value = 990
$ws.Range('AD2').Value = 990
$ws.Range('AF2').Value = 1000
$ws.Range('AG2').Value = 990
$ws.Range('AJ2').Value = 1000
$ws.Range('AN2').Value = 1000

# --- Row 3 updates ---
$ws.Range('C3').Value = '22:00:00'
$ws.Range('D3').Value = 'Alebrijes de Oaxaca'
$ws.Range('E3').Value = 'Chapulineros de Oaxaca'
$ws.Range('F3').Value = 1.34
$ws.Range('G3').Value = 1.42
$ws.Range('H3').Value = 9.199999999999999
$ws.Range('I3').Value = 12
$ws.Range('J3').Value = 5.3
$ws.Range('K3').Value = 6.6
$ws.Range('L3').Value = 1.27
$ws.Range('M3').Value = 1.03
$ws.Range('N3').Value = 5.5
$ws.Range('O3').Value = 1.18
$ws.Range('P3').Value = 2.52
$ws.Range('Q3').Value = 1.56
$ws.Range('R3').Value = 1.59
$ws.Range('S3').Value = 2.44
$ws.Range('T3').Value = 1.86
$ws.Range('U3').Value = 1.89
$ws.Range('V3').Value = 1.09
$ws.Range('W3').Value = 3.35
$ws.Range('X3').Value = 27
$ws.Range('Y3').Value = 46
$ws.Range('Z3').Value = 120
$ws.Range('AA3').Value = 1000
$ws.Range('AB3').Value = 12.5
$ws.Range('AC3').Value = 16.5
$ws.Range('AD3').Value = 44
$ws.Range('AE3').Value = 1000
$ws.Range('AF3').Value = 10
$ws.Range('AG3').Value = 13
$ws.Range('AH3').Value = 32
$ws.Range('AI3').Value = 140
$ws.Range('AJ3').Value = 14
$ws.Range('AK3').Value = 17.5
$ws.Range('AL3').Value = 42
$ws.Range('AM3').Value = 1000
$ws.Range('AN3').Value = 5.1
$ws.Range('AO3').Value = 1000

